$d = $word.ActiveDocument

$sep = "--------------"
$br = [char]11

$quotes = @(
    "I emphasize the importance of screening for gambling disorder in all patients, not just those presenting with substance use issues. Incorporating simple questions about gambling behavior into routine evaluations can help identify this hidden addiction early on.",
    "When gathering a gambling history, I recommend asking patients about their relationship with money and what draws them to gamble. Understanding their motivations can reveal underlying emotional issues, such as loneliness or the need for escape.",
    "For treatment, I often use naltrexone for patients who exhibit strong cravings and urges related to gambling. It’s crucial to address co-occurring disorders, such as ADHD or depression, as treating these can also alleviate gambling issues.",
    "I advocate for a biopsychosocial approach in treatment, focusing first on improving overall health and well-being. By addressing lifestyle factors like sleep, nutrition, and stress management, we can create a foundation that supports recovery from gambling disorder.",
    "It’s essential to recognize that Gamblers Anonymous is only one part of a comprehensive treatment plan. Combining support groups with psychotherapy and medication can lead to better outcomes for patients struggling with gambling addiction.",
    "I advise primary care providers and general mental health clinicians to develop a network of resources for gambling treatment. Knowing where to refer patients for specialized care can significantly enhance their recovery journey.",
    "I’ve found that online resources and telehealth options have expanded access to treatment for gambling disorders. Programs like Kindbridge offer valuable support for patients who may not have local resources available.",
    "For family members of those affected by gambling disorder, I encourage them to seek support through organizations like Gam-Anon. They provide essential resources to help families cope with the impact of gambling addiction."
)

$parts = @()
for ($i = 0; $i -lt $quotes.Length; $i++) {
    $parts += $quotes[$i]
    if ($i -lt $quotes.Length - 1) {
        $parts += $sep
    }
}

$newText = [string]::Join($br, $parts)

$p = $d.Paragraphs.Item(1)
$r = $p.Range
$r.Text = $newText
